# chore: update Sheets via scheduled runner
# Refreshes market-price derived columns (currentAveragePrice* / LevePrice* /
# LeveProfit*) for a handful of Leve rows across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 431.25
$ws.Range("I18").Value = 431.25
$ws.Range("K18").Value = 431.25
$ws.Range("M18").Value = -147.25

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws.Range("H74").Value = 8650
$ws.Range("I74").Value = 8650
$ws.Range("K74").Value = 8650
$ws.Range("M74").Value = -7714

$ws.Range("H76").Value = 455
$ws.Range("I76").Value = 455
$ws.Range("K76").Value = 455
$ws.Range("M76").Value = -140

$ws.Range("H77").Value = 8650
$ws.Range("I77").Value = 8650
$ws.Range("K77").Value = 43250
$ws.Range("M77").Value = -38570

$ws.Range("H79").Value = 455
$ws.Range("I79").Value = 455
$ws.Range("K79").Value = 455
$ws.Range("M79").Value = 637

$ws.Range("H116").Value = 8824.666999999999
$ws.Range("I116").Value = 8249.5
$ws.Range("K116").Value = 8249.5
$ws.Range("M116").Value = -4807.5

$ws.Range("H132").Value = 1206.6957
$ws.Range("I132").Value = 1328.3889
$ws.Range("J132").Value = 768.6
$ws.Range("K132").Value = 3985.1667
$ws.Range("L132").Value = 2305.8
$ws.Range("M132").Value = -1455.1667
$ws.Range("N132").Value = -7365.8

$ws.Range("H138").Value = 4643.9165
$ws.Range("I138").Value = 1578.8
$ws.Range("J138").Value = 5822.8076
$ws.Range("K138").Value = 4736.4
$ws.Range("L138").Value = 17468.4228
$ws.Range("M138").Value = 403.6000000000004
$ws.Range("N138").Value = -27748.4228

$ws.Range("H141").Value = 2749.75
$ws.Range("I141").Value = 999
$ws.Range("K141").Value = 2997
$ws.Range("M141").Value = 2183

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14965.833
$ws.Range("I32").Value = 14965.833
$ws.Range("K32").Value = 14965.833
$ws.Range("M32").Value = -14678.833

$ws.Range("H61").Value = 3888.6667
$ws.Range("I61").Value = 3888.6667
$ws.Range("K61").Value = 3888.6667
$ws.Range("M61").Value = -3676.6667

$ws.Range("H63").Value = 31695.555
$ws.Range("I63").Value = 57184.668
$ws.Range("K63").Value = 57184.668
$ws.Range("M63").Value = -56498.668

$ws.Range("H66").Value = 31695.555
$ws.Range("I66").Value = 57184.668
$ws.Range("K66").Value = 285923.34
$ws.Range("M66").Value = -282491.34

$ws.Range("H74").Value = 10674.286
$ws.Range("I74").Value = 11804.046
$ws.Range("J74").Value = 6531.8335
$ws.Range("K74").Value = 11804.046
$ws.Range("L74").Value = 6531.8335
$ws.Range("M74").Value = -10930.046
$ws.Range("N74").Value = -8279.833500000001

$ws.Range("H77").Value = 10674.286
$ws.Range("I77").Value = 11804.046
$ws.Range("J77").Value = 6531.8335
$ws.Range("K77").Value = 59020.23
$ws.Range("L77").Value = 32659.1675
$ws.Range("M77").Value = -54652.23
$ws.Range("N77").Value = -41395.1675

$ws.Range("H88").Value = 3116.6667
$ws.Range("J88").Value = 3116.6667
$ws.Range("L88").Value = 3116.6667
$ws.Range("N88").Value = -3928.6667

$ws.Range("H91").Value = 3116.6667
$ws.Range("J91").Value = 3116.6667
$ws.Range("L91").Value = 3116.6667
$ws.Range("N91").Value = -5924.6667

$ws.Range("H122").Value = 2823.0715
$ws.Range("J122").Value = 2249.5
$ws.Range("L122").Value = 6748.5
$ws.Range("N122").Value = -11648.5

$ws.Range("H132").Value = 4318.4287
$ws.Range("I132").Value = 3411.3333
$ws.Range("J132").Value = 4998.75
$ws.Range("K132").Value = 10233.9999
$ws.Range("L132").Value = 14996.25
$ws.Range("M132").Value = -7703.999899999999
$ws.Range("N132").Value = -20056.25

$ws.Range("H136").Value = 3888.6667
$ws.Range("I136").Value = 3888.6667
$ws.Range("K136").Value = 11666.0001
$ws.Range("M136").Value = -9116.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6496.6665
$ws.Range("I86").Value = 1996.25
$ws.Range("K86").Value = 1996.25
$ws.Range("M86").Value = -873.25

$ws.Range("H89").Value = 6496.6665
$ws.Range("I89").Value = 1996.25
$ws.Range("K89").Value = 9981.25
$ws.Range("M89").Value = -4365.25

$ws.Range("H134").Value = 5474.1
$ws.Range("I134").Value = 3304.5557
$ws.Range("K134").Value = 9913.667099999999
$ws.Range("M134").Value = -7378.667099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1602.0769
$ws.Range("I31").Value = 1462.8
$ws.Range("K31").Value = 1462.8
$ws.Range("M31").Value = -1167.8

$ws.Range("H34").Value = 1602.0769
$ws.Range("I34").Value = 1462.8
$ws.Range("K34").Value = 1462.8
$ws.Range("M34").Value = -1260.8

$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 1057.3
$ws.Range("I107").Value = 1146.75
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 1146.75
$ws.Range("L107").Value = 699.5
$ws.Range("M107").Value = 773.25
$ws.Range("N107").Value = -4539.5

$ws.Range("H132").Value = 2734.125
$ws.Range("I132").Value = 2227.8333
$ws.Range("K132").Value = 6683.499899999999
$ws.Range("M132").Value = -4153.499899999999

$ws.Range("H134").Value = 5217.727
$ws.Range("I134").Value = 5217.727
$ws.Range("K134").Value = 15653.181
$ws.Range("M134").Value = -13118.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 716642.7
$ws.Range("I4").Value = 2500499.5
$ws.Range("J4").Value = 3100
$ws.Range("K4").Value = 7501498.5
$ws.Range("L4").Value = 9300
$ws.Range("M4").Value = -7501386.5
$ws.Range("N4").Value = -9524

$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3338

$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3204

$ws.Range("H113").Value = 2251
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2251
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6753
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11093

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6131.25
$ws.Range("I80").Value = 2782
$ws.Range("J80").Value = 11713.333
$ws.Range("K80").Value = 2782
$ws.Range("L80").Value = 11713.333
$ws.Range("M80").Value = -1784
$ws.Range("N80").Value = -13709.333

$ws.Range("H83").Value = 6131.25
$ws.Range("I83").Value = 2782
$ws.Range("J83").Value = 11713.333
$ws.Range("K83").Value = 13910
$ws.Range("L83").Value = 58566.665
$ws.Range("M83").Value = -8918
$ws.Range("N83").Value = -68550.66500000001

$ws.Range("H132").Value = 1879.1333
$ws.Range("I132").Value = 1399.1538
$ws.Range("K132").Value = 4197.4614
$ws.Range("M132").Value = -1667.4614

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4666
$ws.Range("I40").Value = 4666
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4666
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4530
$ws.Range("N40").ClearContents()

$ws.Range("H55").Value = 4160
$ws.Range("I55").Value = 5150
$ws.Range("K55").Value = 5150
$ws.Range("M55").Value = -4977

$ws.Range("H132").Value = 4486.375
$ws.Range("I132").Value = 3578.6
$ws.Range("K132").Value = 10735.8
$ws.Range("M132").Value = -8205.799999999999

$ws.Range("H136").Value = 7699.2
$ws.Range("I136").Value = 7699.2
$ws.Range("K136").Value = 23097.6
$ws.Range("M136").Value = -20547.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1845.0769
$ws.Range("I100").Value = 1362.4286
$ws.Range("J100").Value = 2408.1667
$ws.Range("K100").Value = 2724.8572
$ws.Range("L100").Value = 4816.3334
$ws.Range("M100").Value = -2183.8572
$ws.Range("N100").Value = -5898.3334

$ws.Range("H132").Value = 2154.0715
$ws.Range("I132").Value = 707
$ws.Range("K132").Value = 2121
$ws.Range("M132").Value = 409
